{"js": "// The document contains a table of \"three-digit \u00d7 one-digit\" multiplication\n// problems (with results), one expression per table cell, e.g. \"813\u00d74=3252\".\n// This edit swaps each old expression for an updated one (new operands and\n// new result), leaving every other part of the document untouched.\n//\n// Each old expression appears exactly once in the document, so a simple\n// text search + replace per pair is sufficient and preserves the existing\n// run formatting (font/size) of the cell it lives in.\nconst replacements = [\n  [\"813\u00d74=3252\", \"536\u00d73=1608\"],\n  [\"682\u00d79=6138\", \"365\u00d72=730\"],\n  [\"124\u00d75=620\", \"565\u00d79=5085\"],\n  [\"897\u00d76=5382\", \"732\u00d73=2196\"],\n  [\"157\u00d76=942\", \"290\u00d72=580\"],\n  [\"717\u00d75=3585\", \"937\u00d75=4685\"],\n  [\"218\u00d79=1962\", \"974\u00d73=2922\"],\n  [\"286\u00d73=858\", \"486\u00d72=972\"],\n  [\"842\u00d72=1684\", \"841\u00d79=7569\"],\n  [\"438\u00d74=1752\", \"428\u00d78=3424\"],\n  [\"877\u00d78=7016\", \"295\u00d74=1180\"],\n  [\"952\u00d79=8568\", \"188\u00d75=940\"],\n  [\"835\u00d76=5010\", \"136\u00d77=952\"],\n  [\"466\u00d73=1398\", \"186\u00d75=930\"],\n  [\"435\u00d78=3480\", \"778\u00d72=1556\"],\n  [\"729\u00d75=3645\", \"856\u00d73=2568\"],\n  [\"741\u00d75=3705\", \"107\u00d77=749\"],\n  [\"592\u00d76=3552\", \"218\u00d72=436\"],\n  [\"599\u00d79=5391\", \"308\u00d76=1848\"],\n  [\"298\u00d77=2086\", \"110\u00d79=990\"],\n  [\"710\u00d77=4970\", \"643\u00d78=5144\"],\n  [\"240\u00d72=480\", \"309\u00d78=2472\"],\n  [\"684\u00d76=4104\", \"334\u00d73=1002\"],\n  [\"151\u00d74=604\", \"601\u00d79=5409\"],\n  [\"143\u00d73=429\", \"676\u00d75=3380\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains a table of \"three-digit x one-digit\" multiplication\n# problems (with results), one expression per table cell, e.g. \"813x4=3252\".\n# This edit swaps each old expression for an updated one (new operands and\n# new result), leaving every other part of the document untouched.\n#\n# Each old expression appears exactly once in the document, so a simple\n# Find/Replace per pair (restricted to that exact text) is sufficient and\n# preserves the existing run formatting (font/size) of the cell it lives in.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"813\u00d74=3252\", \"536\u00d73=1608\"),\n    @(\"682\u00d79=6138\", \"365\u00d72=730\"),\n    @(\"124\u00d75=620\", \"565\u00d79=5085\"),\n    @(\"897\u00d76=5382\", \"732\u00d73=2196\"),\n    @(\"157\u00d76=942\", \"290\u00d72=580\"),\n    @(\"717\u00d75=3585\", \"937\u00d75=4685\"),\n    @(\"218\u00d79=1962\", \"974\u00d73=2922\"),\n    @(\"286\u00d73=858\", \"486\u00d72=972\"),\n    @(\"842\u00d72=1684\", \"841\u00d79=7569\"),\n    @(\"438\u00d74=1752\", \"428\u00d78=3424\"),\n    @(\"877\u00d78=7016\", \"295\u00d74=1180\"),\n    @(\"952\u00d79=8568\", \"188\u00d75=940\"),\n    @(\"835\u00d76=5010\", \"136\u00d77=952\"),\n    @(\"466\u00d73=1398\", \"186\u00d75=930\"),\n    @(\"435\u00d78=3480\", \"778\u00d72=1556\"),\n    @(\"729\u00d75=3645\", \"856\u00d73=2568\"),\n    @(\"741\u00d75=3705\", \"107\u00d77=749\"),\n    @(\"592\u00d76=3552\", \"218\u00d72=436\"),\n    @(\"599\u00d79=5391\", \"308\u00d76=1848\"),\n    @(\"298\u00d77=2086\", \"110\u00d79=990\"),\n    @(\"710\u00d77=4970\", \"643\u00d78=5144\"),\n    @(\"240\u00d72=480\", \"309\u00d78=2472\"),\n    @(\"684\u00d76=4104\", \"334\u00d73=1002\"),\n    @(\"151\u00d74=604\", \"601\u00d79=5409\"),\n    @(\"143\u00d73=429\", \"676\u00d75=3380\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
